# Update cryptocurrency price (D) and volume change (E) figures
# per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain decimal number need
# to be pre-formatted as Text so Excel does not silently convert the
# price string into a numeric value (which would drop formatting such as
# trailing zeros, e.g. '1.00' or '0.530').
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "46.780.93"
$ws.Range("E2").Value = "  +6.46%  "
$ws.Range("D3").Value = "2.310.03"
$ws.Range("E3").Value = "  +5.10%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "299.55"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "102.25"
$ws.Range("E6").Value = "  +15.23%  "
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +10.08%  "
$ws.Range("D10").Value = "36.85"
$ws.Range("E10").Value = "  +14.99%  "
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  +4.52%  "
$ws.Range("E12").Value = "  +9.11%  "
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "2.659.97"
$ws.Range("E14").Value = "  +5.10%  "
$ws.Range("D15").Value = "2.307.69"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "13.99"
$ws.Range("E16").Value = "  +7.26%  "
$ws.Range("D17").Value = "0.821"
$ws.Range("E17").Value = "  +6.54%  "
$ws.Range("D18").Value = "46.752.00"
$ws.Range("E18").Value = "  +7.29%  "
$ws.Range("E19").Value = "  +23.01%  "
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  +6.71%  "
$ws.Range("E21").Value = "  +5.49%  "
$ws.Range("D22").Value = "67.13"
$ws.Range("E22").Value = "  +6.68%  "
$ws.Range("D23").Value = "248.31"
$ws.Range("E23").Value = "  +7.66%  "
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  +6.51%  "
$ws.Range("E25").Value = "  +8.05%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "42.98"
$ws.Range("E27").Value = "  +19.61%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").Value = "9.93"
$ws.Range("E29").Value = "  +7.93%  "
$ws.Range("E30").Value = "  +5.72%  "
$ws.Range("D31").Value = "5.77"
$ws.Range("E31").Value = "  +8.85%  "
$ws.Range("D32").Value = "147.18"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "0.0801"
$ws.Range("E33").Value = "  +9.16%  "
$ws.Range("D34").Value = "2.63"
$ws.Range("E34").Value = "  +4.65%  "
$ws.Range("E35").Value = "  +8.45%  "
$ws.Range("E36").Value = "  +9.80%  "
$ws.Range("D37").Value = "0.120"
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("E38").Value = "  +8.90%  "
$ws.Range("D39").Value = "15.74"
$ws.Range("E39").Value = "  +20.91%  "
$ws.Range("D40").Value = "4.03"
$ws.Range("E40").Value = "  +14.89%  "
$ws.Range("D41").Value = "3.44"
$ws.Range("E41").Value = "  +11.37%  "
$ws.Range("E42").Value = "  +9.23%  "
$ws.Range("D43").Value = "2.02"
$ws.Range("E43").Value = "  +23.70%  "
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "1.843.45"
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").Value = "89.62"
$ws.Range("E46").Value = "  +23.89%  "
$ws.Range("E47").Value = "  +14.61%  "
$ws.Range("D48").Value = "75.15"
$ws.Range("E48").Value = "  +16.04%  "
$ws.Range("D49").Value = "4.95"
$ws.Range("E49").Value = "  +9.64%  "
$ws.Range("D50").Value = "97.55"
$ws.Range("E50").Value = "  +6.35%  "
$ws.Range("D51").Value = "54.56"
$ws.Range("E51").Value = "  +11.45%  "
